$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: Row number, new Price (column D) text, new Volume(1h) (column E) text.
# A $null for D means that cell's price text is unchanged in this update.
$updates = @(
    @{ Row = 2; D = '29.195.40'; E = '  +0.34%  ' },
    @{ Row = 3; D = '1.833.54'; E = '  -0.17%  ' },
    @{ Row = 4; D = $null; E = '  +0.02%  ' },
    @{ Row = 5; D = '242.28'; E = '  +0.74%  ' },
    @{ Row = 6; D = '0.6594'; E = '  -2.64%  ' },
    @{ Row = 7; D = $null; E = '  +0.03%  ' },
    @{ Row = 8; D = '0.07407'; E = '  -0.64%  ' },
    @{ Row = 9; D = '0.2934'; E = '  -1.63%  ' },
    @{ Row = 10; D = '22.87'; E = '  -0.85%  ' },
    @{ Row = 11; D = '0.07774'; E = '  +1.45%  ' },
    @{ Row = 12; D = '1.815.08'; E = '  -1.17%  ' },
    @{ Row = 13; D = '5.000'; E = '  -0.43%  ' },
    @{ Row = 14; D = '0.6665'; E = '  -1.40%  ' },
    @{ Row = 15; D = '82.83'; E = '  -3.98%  ' },
    @{ Row = 16; D = '6.108'; E = '  -0.87%  ' },
    @{ Row = 17; D = '0.000008405'; E = '  +1.79%  ' },
    @{ Row = 18; D = '29.149.82'; E = '  +0.20%  ' },
    @{ Row = 19; D = '2.054.40'; E = '  -1.24%  ' },
    @{ Row = 20; D = '227.38'; E = '  -0.10%  ' },
    @{ Row = 21; D = '12.47'; E = '  -0.13%  ' },
    @{ Row = 22; D = '1.000'; E = '  +0.14%  ' },
    @{ Row = 23; D = '7.128'; E = '  -2.93%  ' },
    @{ Row = 24; D = $null; E = '  +0.05%  ' },
    @{ Row = 25; D = '158.80'; E = '  -1.13%  ' },
    @{ Row = 26; D = '8.610'; E = '  -1.07%  ' },
    @{ Row = 27; D = '0.1393'; E = '  -2.68%  ' },
    @{ Row = 28; D = '17.95'; E = '  -0.32%  ' },
    @{ Row = 29; D = $null; E = '  +0.87%  ' },
    @{ Row = 30; D = '4.114'; E = '  -3.08%  ' },
    @{ Row = 31; D = '4.048'; E = '  -2.07%  ' },
    @{ Row = 32; D = '1.194'; E = '  -0.31%  ' },
    @{ Row = 33; D = '0.05268'; E = '  -3.18%  ' },
    @{ Row = 34; D = '1.865'; E = '  +0.29%  ' },
    @{ Row = 35; D = '0.7425'; E = '  -0.72%  ' },
    @{ Row = 36; D = $null; E = '  +1.19%  ' },
    @{ Row = 37; D = '2.653'; E = '  -1.08%  ' },
    @{ Row = 38; D = '1.306.00'; E = '  +0.19%  ' },
    @{ Row = 39; D = '0.01792'; E = '  -1.05%  ' },
    @{ Row = 40; D = '2.736'; E = '  +1.04%  ' },
    @{ Row = 41; D = '0.9293'; E = '  -0.59%  ' },
    @{ Row = 42; D = '5.929'; E = '  -2.68%  ' },
    @{ Row = 43; D = '0.08433'; E = '  +5.01%  ' },
    @{ Row = 44; D = '0.9998'; E = '  +0.07%  ' },
    @{ Row = 45; D = '102.55'; E = '  -1.84%  ' },
    @{ Row = 46; D = '1.959.24'; E = '  -0.84%  ' },
    @{ Row = 47; D = '0.5142'; E = '  -0.61%  ' },
    @{ Row = 48; D = $null; E = '  -1.12%  ' },
    @{ Row = 49; D = '1.750'; E = '  -0.58%  ' },
    @{ Row = 50; D = '62.92'; E = '  -1.13%  ' },
    @{ Row = 51; D = '0.05866'; E = '  -1.09%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        # Values like "5.000" or "29.195.40" must stay text (not be coerced into
        # numbers - which would also drop significant trailing zeros, e.g.
        # "5.000" -> 5), so force the column-D cell to a text format before
        # writing it, matching the text already used by this sheet. Restore the
        # cell's style back to Normal afterwards so we don't leave a stray
        # "@"-formatted look on the cell.
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }

    if ($null -ne $u.E) {
        $eCell = $ws.Cells.Item($row, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $u.E
        $eCell.Style = "Normal"
    }
}
